# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.133.36"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.842.53"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'241.08"
$ws.Range("E5").Value = "  -2.03%  "

$ws.Range("D6").Value = "'0.6863"
$ws.Range("E6").Value = "  -1.74%  "

$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -1.33%  "

$ws.Range("D9").Value = "'0.07473"
$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("D10").Value = "'23.13"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("D11").Value = "'0.07666"
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").Value = "1.842.51"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").Value = "'0.6843"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "'87.59"
$ws.Range("E15").Value = "  -5.84%  "

$ws.Range("E16").Value = "  -7.04%  "

$ws.Range("D17").Value = "29.130.57"
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").Value = "'0.000008161"
$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("D19").Value = "2.079.92"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "'228.46"
$ws.Range("E20").Value = "  -5.34%  "

$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'7.406"
$ws.Range("E23").Value = "  -1.47%  "

$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  -3.66%  "

$ws.Range("D26").Value = "'159.95"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").Value = "'8.761"
$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("D30").Value = "'4.266"
$ws.Range("E30").Value = "  +0.90%  "

$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").Value = "'0.05211"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").Value = "'0.7658"
$ws.Range("E34").Value = "  -3.66%  "

$ws.Range("D35").Value = "'1.851"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").Value = "'1.136"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").Value = "1.315.53"
$ws.Range("E38").Value = "  -0.16%  "

$ws.Range("D39").Value = "'0.01838"
$ws.Range("E39").Value = "  -1.92%  "

$ws.Range("D40").Value = "'2.726"
$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("D41").Value = "'0.9334"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("D42").Value = "'104.84"
$ws.Range("E42").Value = "  -2.06%  "

$ws.Range("D43").Value = "'5.786"
$ws.Range("E43").Value = "  -3.69%  "

$ws.Range("D44").Value = "'0.9999"

$ws.Range("D45").Value = "1.982.46"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").Value = "'64.87"
$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.565"
$ws.Range("E48").Value = "  -1.52%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("E49").Value = "  -2.12%  "

$ws.Range("E50").Value = "  +0.59%  "

$ws.Range("D51").Value = "'0.07437"
$ws.Range("E51").Value = "  +17.79%  "
